$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells keep their original Text format so Excel
# does not auto-convert numeric-looking strings (e.g. "0.9998",
# "27.81", "0.000008078") into real numbers / scientific notation.
$cells = @(
    "D2", "E2", "D3", "E3", "E4", "D5", "E5", "D6", "E6", "D7", "E8", "D9", "D10", "E10", "D11", "E11", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "E16", "D17", "E17", "D18", "E18", "D19", "E19", "D20", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "D26", "E26", "E27", "D28", "E28", "D29", "E29", "E30", "D31", "E31", "E32", "E33", "D34", "E34", "E35", "D36", "E36", "D37", "E37", "D38", "E38", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "D47", "E47", "D48", "E48", "D49", "E49", "D50", "E50", "D51", "E51"
)
foreach ($c in $cells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = "30.292.21"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").Value = "1.928.16"
$ws.Range("E3").Value = "  -0.76%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "249.06"
$ws.Range("E5").Value = "  -0.81%  "
$ws.Range("D6").Value = "0.7177"
$ws.Range("E6").Value = "  -0.72%  "
$ws.Range("D7").Value = "0.9998"
$ws.Range("E8").Value = "  -4.97%  "
$ws.Range("D9").Value = "27.81"
$ws.Range("D10").Value = "0.07074"
$ws.Range("E10").Value = "  -4.53%  "
$ws.Range("D11").Value = "0.7890"
$ws.Range("E11").Value = "  -3.57%  "
$ws.Range("E12").Value = "  -1.95%  "
$ws.Range("D13").Value = "1.929.13"
$ws.Range("E13").Value = "  -0.65%  "
$ws.Range("D14").Value = "5.375"
$ws.Range("E14").Value = "  -2.82%  "
$ws.Range("D15").Value = "94.70"
$ws.Range("E15").Value = "  -0.83%  "
$ws.Range("E16").Value = "  -1.79%  "
$ws.Range("D17").Value = "30.292.99"
$ws.Range("E17").Value = "  -0.22%  "
$ws.Range("D18").Value = "256.91"
$ws.Range("E18").Value = "  +0.81%  "
$ws.Range("D19").Value = "0.000008078"
$ws.Range("E19").Value = "  -2.65%  "
$ws.Range("D20").Value = "5.752"
$ws.Range("E20").Value = "  -2.47%  "
$ws.Range("D21").Value = "2.184.13"
$ws.Range("E21").Value = "  -0.45%  "
$ws.Range("D22").Value = "0.9999"
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("D24").Value = "6.820"
$ws.Range("E24").Value = "  -2.13%  "
$ws.Range("D25").Value = "9.533"
$ws.Range("E25").Value = "  -3.55%  "
$ws.Range("D26").Value = "164.39"
$ws.Range("E26").Value = "  +2.43%  "
$ws.Range("E27").Value = "  -1.92%  "
$ws.Range("D28").Value = "2.267"
$ws.Range("E28").Value = "  -7.95%  "
$ws.Range("D29").Value = "0.1270"
$ws.Range("E29").Value = "  -4.00%  "
$ws.Range("E30").Value = "  +0.88%  "
$ws.Range("D31").Value = "1.529"
$ws.Range("E31").Value = "  -2.63%  "
$ws.Range("E32").Value = "  -2.13%  "
$ws.Range("E33").Value = "  -2.68%  "
$ws.Range("D34").Value = "0.05127"
$ws.Range("E34").Value = "  -2.65%  "
$ws.Range("E35").Value = "  -1.03%  "
$ws.Range("D36").Value = "0.7449"
$ws.Range("E36").Value = "  -1.59%  "
$ws.Range("D37").Value = "2.766"
$ws.Range("E37").Value = "  +1.02%  "
$ws.Range("D38").Value = "0.01981"
$ws.Range("E38").Value = "  -0.57%  "
$ws.Range("D40").Value = "78.08"
$ws.Range("E40").Value = "  -4.79%  "
$ws.Range("D41").Value = "6.371"
$ws.Range("E41").Value = "  -2.93%  "
$ws.Range("D42").Value = "0.4504"
$ws.Range("E42").Value = "  -2.03%  "
$ws.Range("D43").Value = "1.995"
$ws.Range("E43").Value = "  -1.83%  "
$ws.Range("D44").Value = "0.8460"
$ws.Range("E44").Value = "  -0.54%  "
$ws.Range("D45").Value = "0.9995"
$ws.Range("E45").Value = "  -0.06%  "
$ws.Range("D46").Value = "100.72"
$ws.Range("E46").Value = "  -2.40%  "
$ws.Range("D47").Value = "9.798"
$ws.Range("E47").Value = "  -1.18%  "
$ws.Range("D48").Value = "7.448"
$ws.Range("E48").Value = "  -0.35%  "
$ws.Range("D49").Value = "36.84"
$ws.Range("E49").Value = "  -0.78%  "
$ws.Range("D50").Value = "951.22"
$ws.Range("E50").Value = "  +7.76%  "
$ws.Range("D51").Value = "0.4211"
$ws.Range("E51").Value = "  -0.06%  "
